# edit.ps1
# Implements commit: "Additional processing of the URL field is empty"
#
# The worksheet "咕咕咕" lists business systems together with a URL and an
# IP column. This change appends two more systems to the bottom of the
# table:
#   - row 11: "哇嘎嘎URL是空的系统" - a system whose URL column is left
#             completely empty (the scenario the commit message refers to)
#   - row 12: "啦啦啦" - a normal system with a URL (and its hyperlink),
#             just like the other rows above it
#
# It also leaves the sheet with the selection the author ended up with
# (cell C6) and keeps the column widths very close to the saved ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 11: a system whose URL field is empty ----------------------------
# Copy the plain (unstyled) format of A10 first so the new cell does not
# pick up an explicit style index, matching how column A looks elsewhere.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "哇嘎嘎URL是空的系统"
# (B11/C11 intentionally left blank - no URL, no IP for this system)

# --- Row 12: a normal system with a URL + hyperlink -----------------------
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = "啦啦啦"

$ws.Range("B12").Value = "https://192.111.23.56"
$ws.Hyperlinks.Add($ws.Range("B12"), "https://192.111.23.56")
# Re-apply the same hyperlink-cell formatting used by the other URL cells
# (e.g. B10) so B12 matches the look of the rest of the table.
$ws.Range("B10").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Column widths (kept close to the saved widths) ------------------------
$ws.Columns.Item(1).ColumnWidth = 19.551339285714285
$ws.Columns.Item(2).ColumnWidth = 33.883370535714285
$ws.Columns.Item(3).ColumnWidth = 25.215401785714285

# --- Leftover UI state: active selection left on C6 ------------------------
$ws.Range("C6").Select() | Out-Null
